$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 141; $r++) {
    $ws.Cells.Item($r, 4).Value = 153
    $ws.Cells.Item($r, 5).Value = 872
}

$ws.Range("F8").Select()
